$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.804.00"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "1.925.33"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("E7").Value = "  -1.93%  "

$ws.Range("D8").Value = "'0.2886"
$ws.Range("E8").Value = "  -2.27%  "

$ws.Range("D9").Value = "'0.06773"
$ws.Range("E9").Value = "  -1.72%  "

$ws.Range("D10").Value = "'19.61"
$ws.Range("E10").Value = "  +1.48%  "

$ws.Range("D11").Value = "'103.89"
$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("D12").Value = "'0.07795"
$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("D13").Value = "1.933.21"
$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("D14").Value = "'5.268"
$ws.Range("E14").Value = "  -1.50%  "

$ws.Range("D15").Value = "'0.6818"
$ws.Range("E15").Value = "  -2.90%  "

$ws.Range("D16").Value = "'292.23"
$ws.Range("E16").Value = "  +7.09%  "

$ws.Range("D17").Value = "30.828.73"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("B18").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C18").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D18").Value = "2.191.64"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007576"
$ws.Range("E19").Value = "  -1.88%  "

$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").Value = "'12.86"
$ws.Range("E21").Value = "  -1.92%  "

$ws.Range("D22").Value = "'5.518"
$ws.Range("E22").Value = "  -3.17%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'6.379"
$ws.Range("E24").Value = "  -2.44%  "

$ws.Range("D25").Value = "'9.539"
$ws.Range("E25").Value = "  -2.84%  "

$ws.Range("D26").Value = "'167.98"
$ws.Range("E26").Value = "  +1.70%  "

$ws.Range("D27").Value = "'19.73"
$ws.Range("E27").Value = "  +0.73%  "

$ws.Range("D28").Value = "'2.115"
$ws.Range("E28").Value = "  -2.15%  "

$ws.Range("D29").Value = "'1.392"
$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("D30").Value = "'0.1007"
$ws.Range("E30").Value = "  -2.71%  "

$ws.Range("D31").Value = "'4.595"
$ws.Range("E31").Value = "  -2.26%  "

$ws.Range("D32").Value = "'1.527"
$ws.Range("E32").Value = "  -1.88%  "

$ws.Range("D33").Value = "'4.322"
$ws.Range("E33").Value = "  -2.49%  "

$ws.Range("D34").Value = "'0.04812"
$ws.Range("E34").Value = "  -1.92%  "

$ws.Range("D35").Value = "'0.7343"
$ws.Range("E35").Value = "  -3.36%  "

$ws.Range("D36").Value = "'1.124"
$ws.Range("E36").Value = "  -2.31%  "

$ws.Range("D37").Value = "'2.720"
$ws.Range("E37").Value = "  -0.46%  "

$ws.Range("D38").Value = "'0.01941"
$ws.Range("E38").Value = "  -3.34%  "

$ws.Range("D39").Value = "'2.633"
$ws.Range("E39").Value = "  -1.54%  "

$ws.Range("D40").Value = "'6.408"
$ws.Range("E40").Value = "  -0.93%  "

$ws.Range("D41").Value = "'75.28"
$ws.Range("E41").Value = "  -4.96%  "

$ws.Range("D42").Value = "'2.012"
$ws.Range("E42").Value = "  -3.12%  "

$ws.Range("D43").Value = "'0.8657"
$ws.Range("E43").Value = "  -3.84%  "

$ws.Range("D44").Value = "'0.4334"
$ws.Range("E44").Value = "  -2.70%  "

$ws.Range("D45").Value = "'105.64"
$ws.Range("E45").Value = "  -2.54%  "

$ws.Range("D46").Value = "'0.9998"
$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("D47").Value = "'7.529"
$ws.Range("E47").Value = "  -4.57%  "

$ws.Range("D48").Value = "'998.45"
$ws.Range("E48").Value = "  +1.02%  "

$ws.Range("D49").Value = "'0.1209"
$ws.Range("E49").Value = "  -3.24%  "

$ws.Range("D50").Value = "'9.050"
$ws.Range("E50").Value = "  -2.37%  "

$ws.Range("D51").Value = "'34.91"
$ws.Range("E51").Value = "  -3.75%  "
